# Restructure the "Company, City, ST" Heading3 lines in the Experience /
# Skills sections into separate Heading3 (company) + Heading4 (location)
# paragraphs, and demote the job-title / skill-subheading lines that used
# to be Heading4 down to Heading5 so the new location line can take their
# old outline slot.
#
# NOTE: this interpreter does not bind PowerShell named (`-Param value`)
# arguments, so every helper function below uses plain positional params.

$d = $word.ActiveDocument

# --- Helper: find the 1-based index of the paragraph whose visible text
#     (paragraph mark stripped) equals $Text.
function Get-ParagraphIndexByText($Text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -eq $Text) {
            return $i
        }
    }
    return -1
}

function Split-CompanyLocation($FullText, $TitleText, $NewTitleBookmark, $LocationText, $LocationBookmark) {
    $idx = Get-ParagraphIndexByText $FullText
    if ($idx -lt 0) {
        throw "Paragraph not found: $FullText"
    }
    $p = $d.Paragraphs($idx)

    # Trim the paragraph text down to just the title (drop ", Location").
    $r = $p.Range
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $TitleText

    # Rename the bookmark wrapping the title: drop whatever bookmark(s)
    # currently sit on this paragraph and add the renamed one back.
    $titleRange = $d.Paragraphs($idx).Range.Duplicate
    $titleRange.MoveEnd(1, -1) | Out-Null
    foreach ($bm in @($d.Bookmarks)) {
        if ($bm.Start -ge $p.Range.Start -and $bm.End -le $p.Range.End) {
            $bm.Delete()
        }
    }
    $d.Bookmarks.Add($NewTitleBookmark, $titleRange) | Out-Null

    # Insert a new Heading4 paragraph right after for the location.
    $insPoint = $d.Paragraphs($idx).Range.Duplicate
    $insPoint.InsertParagraphAfter() | Out-Null
    $locPara = $d.Paragraphs($idx + 1)
    $locPara.Style = "Heading4"
    $locRange = $locPara.Range
    $locRange.MoveEnd(1, -1) | Out-Null
    $locRange.Text = $LocationText
    $d.Bookmarks.Add($LocationBookmark, $locRange) | Out-Null

    # The paragraph that used to follow the combined heading (old Heading4
    # job-title / subheading line) now needs to move down one outline
    # level, from Heading4 to Heading5. It kept its own bookmark name.
    $nextPara = $d.Paragraphs($idx + 2)
    $nextPara.Style = "Heading5"
}

Split-CompanyLocation "Black Beard Labs (formerly Veteran Studio), Waynesboro, VA" `
    "Black Beard Labs (formerly Veteran Studio)" "black-beard-labs-formerly-veteran-studio" `
    "Waynesboro, VA" "waynesboro-va"

Split-CompanyLocation "Innovative Refrigeration Systems, Lyndhurst, VA" `
    "Innovative Refrigeration Systems" "innovative-refrigeration-systems" `
    "Lyndhurst, VA" "lyndhurst-va"

Split-CompanyLocation "Pacific Crest Trail, CA-OR-WA, USA" `
    "Pacific Crest Trail" "pacific-crest-trail" `
    "CA-OR-WA, USA" "ca-or-wa-usa"

Split-CompanyLocation "MidwayUSA, Inc., Columbia, MO" `
    "MidwayUSA, Inc." "midwayusa-inc." `
    "Columbia, MO" "columbia-mo"

Split-CompanyLocation "Zeta-Meter, Inc., Staunton, VA" `
    "Zeta-Meter, Inc." "zeta-meter-inc." `
    "Staunton, VA" "staunton-va"

# Standalone Heading4 -> Heading5 demotions (no location split involved).
function Set-ParagraphStyleByText($Text, $StyleName) {
    $idx = Get-ParagraphIndexByText $Text
    if ($idx -lt 0) {
        throw "Paragraph not found: $Text"
    }
    $d.Paragraphs($idx).Style = $StyleName
}

Set-ParagraphStyleByText "Application Developer – 2010-2013" "Heading5"
Set-ParagraphStyleByText "Polyglot" "Heading5"
Set-ParagraphStyleByText "Specific Solutions and Software Packages" "Heading5"
Set-ParagraphStyleByText "Areas of Focus" "Heading5"

Write-Output "done"
